# Add team Wins/Losses/Ties record columns (AD:AF) to the player data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new column headers
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the existing header formatting (bold, bordered, centered) used by
# the other header cells, e.g. AC1, by copying its format onto the new ones.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Data rows 2-68: every row gets the same team record values.
for ($r = 2; $r -le 68; $r++) {
    $ws.Cells.Item($r, 30).Value = 73  # AD: Wins
    $ws.Cells.Item($r, 31).Value = 89  # AE: Losses
    $ws.Cells.Item($r, 32).Value = 0   # AF: Ties
}
